$d = $word.ActiveDocument

# Locate the paragraph ending in "LOQ4233: Gestão de Negócios (Requisito fraco)".
# Deletion must start right after this paragraph's own ending mark, so grow
# the found range by one character (its paragraph mark) before reading .End.
$startRange = $d.Content
$startRange.Find.Execute("LOQ4233: Gestão de Negócios (Requisito fraco)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startRange.MoveEnd(1, 1)
$delStart = $startRange.End

# Locate the paragraph containing the trailing copyright/footer text; the
# deletion should end right after this paragraph's own ending mark too, so it
# gets removed along with the empty paragraph and the "Ver no Jupiter..." one.
$endRange = $d.Content
$endRange.Find.Execute("© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endRange.MoveEnd(1, 1)
$delEnd = $endRange.End

# Remove the empty paragraph, the "Ver no Jupiter..." paragraph, and the
# "© 2020 ..." paragraph in one shot.
$d.Range($delStart, $delEnd).Delete()
